# Fruta / hortaliza, semanal
# New weekly records were inserted at the top of the dated block
# (rows 108-109), pushing all existing records from row 108 down by
# two rows (to rows 110-226). The two new rows reuse the
# variety/quality/price/weight-price fields of the rows they
# displaced, but carry a new date and new volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 108; this pushes the
# existing rows 108:224 down to 110:226 (and grows the used range to
# A1:T226), preserving all of their values/formatting untouched.
$ws.Rows("108:109").Insert()

# The two new blank rows just inherited the date number format from
# the row above. Populate them with a copy of what is now directly
# below them (the records that used to be at 108/109), then overwrite
# just the date (D) and volume (M) fields with the new values.
$ws.Range("A110:T110").Copy()
$ws.Range("A108:T108").PasteSpecial()

$ws.Range("A111:T111").Copy()
$ws.Range("A109:T109").PasteSpecial()

$ws.Range("D108").Value = 44601
$ws.Range("M108").Value = 60

$ws.Range("D109").Value = 44601
$ws.Range("M109").Value = 30
